$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A171").Value = "2023-12-11 08:28:50"
$ws.Range("B171").Value = 0.0012

$ws.Range("A172").Value = "2023-12-11 08:29:30"
$ws.Range("B172").Value = 0.003200000000000001

$ws.Range("A173").Value = "2023-12-11 08:30:11"
$ws.Range("B173").Value = 0.003
